# Update (Removed Auto Arima)
# Recomputed weekly forecast figures on the "Forecast Comparison" sheet now
# that the Auto-ARIMA model has been dropped from the ensemble, and rolled
# those new numbers up into the "Summary" sheet's headline statistics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Forecast Comparison" sheet: Prophet / Amazon Mean / P70 / P80 / P90
#    forecasts per week (columns C:G, rows 2-17) with the Auto Arima model
#    removed from the blend.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Forecast Comparison")

# row -> (Prophet Forecast, Amazon Mean Forecast, Amazon P70, Amazon P80, Amazon P90)
$newForecast = [ordered]@{
    2  = @(19, 19, 23, 30, 41)
    3  = @(17, 15, 18, 26, 39)
    4  = @(15, 11, 13, 19, 28)
    5  = @(13, 10, 12, 17, 26)
    6  = @(12, 10, 11, 16, 25)
    7  = @(12, 10, 11, 16, 25)
    8  = @(12, 10, 11, 17, 26)
    9  = @(10,  9, 10, 15, 26)
    10 = @(11,  9, 10, 16, 25)
    11 = @(12,  9, 10, 16, 25)
    12 = @(12,  9, 10, 15, 25)
    13 = @(12, 10, 11, 17, 27)
    14 = @(12, 10, 10, 16, 26)
    15 = @(12,  9,  9, 15, 26)
    16 = @(12, 10, 10, 16, 27)
    17 = @(11,  9,  9, 15, 26)
}

foreach ($row in $newForecast.Keys) {
    $vals = $newForecast[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]   # C - Prophet Forecast
    $ws.Cells.Item($row, 4).Value = $vals[1]   # D - Amazon Mean Forecast
    $ws.Cells.Item($row, 5).Value = $vals[2]   # E - Amazon P70 Forecast
    $ws.Cells.Item($row, 6).Value = $vals[3]   # F - Amazon P80 Forecast
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - Amazon P90 Forecast
}

# ---------------------------------------------------------------------------
# 2) "Summary" sheet: headline forecast totals, recomputed from the updated
#    Prophet Forecast column above. The "Value" column stores its entries
#    as text (matching every other row on this sheet), so re-enter each
#    number with a leading apostrophe to keep it a text value rather than
#    letting Excel auto-convert it to a number.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summaryUpdates = [ordered]@{
    9  = "204"   # Total Forecast (16 Weeks)
    10 = "110"   # Total Forecast (8 Weeks)
    11 = "64"    # Total Forecast (4 Weeks)
    12 = "19"    # Max Forecast
    14 = "10"    # Min Forecast
}

foreach ($row in $summaryUpdates.Keys) {
    $summary.Cells.Item($row, 2).Value = "'" + $summaryUpdates[$row]
}
